$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "...nuestro banco de datos de la estación meteorológica" ->
#           "...nuestra base de datos de la estación meteorológica"
# In the target the sentence stays split across three runs, matching what
# Word does when a user selects just the "o banco " substring and retypes
# "a base " (prefix run / retyped run / unchanged suffix run). We reproduce
# that precisely: find the narrow "o banco " span, replace its text, and
# nudge a character property on it (set then reset) so the run boundary
# introduced by the edit is preserved instead of being re-merged with its
# neighbours.
# ---------------------------------------------------------------------------
$mid = $d.Content
$mid.Find.Execute("o banco ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$mid.Text = "a base "
$mid.Bold = 1
$mid.Bold = 0

# ---------------------------------------------------------------------------
# Change 2: the cached TIME field result changes from "1 de noviembre de
# 2024" to "24 de febrero de 2025".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("1 de noviembre de 2024", $true, $false, $false, $false, $false, $true, 1, $false, "24 de febrero de 2025", 2)
